$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last two data rows (rows 4 and 5) — the TPM update collapses
# the table from 4 pairs down to 2.
$ws.Rows("4:5").Delete()

# Row 2: ECs -> Resolving-Mac (target cluster) with refreshed TPM stats.
$ws.Range("D2").Value2 = "Resolving-Mac"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.02345233333333334
$ws.Range("H2").Value2 = 0.070357
$ws.Range("I2").Value2 = 0.0002537772683371841
$ws.Range("J2").Value2 = 0.0002537772683371841
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 5.273410666666667
$ws.Range("N2").Value2 = 15.820232
$ws.Range("O2").Value2 = 1
$ws.Range("P2").Value2 = 1
$ws.Range("Q2").Value2 = 0.1236737847582222
$ws.Range("R2").Value2 = 1.113064062824
$ws.Range("S2").Value2 = 0.0002537772683371841
$ws.Range("T2").Value2 = 0.0002537772683371841

# Row 3: Sending cluster ECs -> Resolving-Mac, with refreshed TPM stats.
$ws.Range("A3").Value2 = "Resolving-Mac"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 92.38960533333334
$ws.Range("H3").Value2 = 277.168816
$ws.Range("I3").Value2 = 0.9997462227316628
$ws.Range("J3").Value2 = 0.9997462227316628
$ws.Range("M3").Value2 = 5.273410666666667
$ws.Range("N3").Value2 = 15.820232
$ws.Range("O3").Value2 = 1
$ws.Range("P3").Value2 = 1
$ws.Range("Q3").Value2 = 487.2083302539236
$ws.Range("R3").Value2 = 4384.874972285312
$ws.Range("S3").Value2 = 0.9997462227316628
$ws.Range("T3").Value2 = 0.9997462227316628
